$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Price" (D) and "Volume(1h)" (E) columns for the rows that changed.
# D-column values are set via a Text number format so digit strings like
# "566.71" or "0.0723" are kept as literal text instead of being parsed as numbers,
# matching the original inline-string cell content; formats are cleared again
# afterwards so no extra cell style is introduced.
$priceUpdates = @{
    2 = '69.904.84'
    3 = '2.457.74'
    5 = '566.71'
    6 = '167.02'
    8 = '0.512'
    9 = '0.176'
    11 = '0.336'
    12 = '4.68'
    14 = '69.760.63'
    15 = '2.910.38'
    16 = '24.04'
    17 = '2.457.49'
    18 = '10.82'
    19 = '343.46'
    20 = '7.13'
    21 = '3.88'
    22 = '2.01'
    24 = '66.59'
    25 = '3.92'
    26 = '2.590.96'
    27 = '8.51'
    28 = '0.991'
    29 = '0.0₃0852'
    30 = '7.32'
    32 = '449.09'
    35 = '161.17'
    39 = '18.15'
    40 = '0.305'
    41 = '1.55'
    42 = '4.48'
    43 = '1.09'
    44 = '2.15'
    45 = '3.40'
    46 = '131.78'
    47 = '0.0723'
    48 = '0.491'
    49 = '0.563'
}

$volumeUpdates = @{
    2 = '  +2.26%  '
    3 = '  +0.49%  '
    4 = '  -0.04%  '
    5 = '  +1.65%  '
    6 = '  +3.01%  '
    7 = '  -0.04%  '
    8 = '  +0.43%  '
    9 = '  +12.92%  '
    10 = '  -1.32%  '
    11 = '  +3.02%  '
    12 = '  -2.61%  '
    13 = '  +8.04%  '
    14 = '  +2.17%  '
    15 = '  +0.16%  '
    16 = '  +3.26%  '
    17 = '  +0.18%  '
    18 = '  +3.74%  '
    19 = '  +2.09%  '
    20 = '  +3.45%  '
    21 = '  +2.56%  '
    22 = '  +6.40%  '
    23 = '  +0.00%  '
    24 = '  -0.09%  '
    25 = '  +6.27%  '
    26 = '  +0.77%  '
    27 = '  +4.35%  '
    28 = '  -0.69%  '
    29 = '  +4.40%  '
    30 = '  +1.79%  '
    31 = '  +9.64%  '
    32 = '  +5.31%  '
    33 = '  +0.12%  '
    34 = '  +1.19%  '
    35 = '  +0.25%  '
    36 = '  +0.43%  '
    37 = '  -0.01%  '
    39 = '  +2.13%  '
    40 = '  +2.97%  '
    41 = '  +5.49%  '
    42 = '  +2.50%  '
    43 = '  +2.71%  '
    44 = '  +6.41%  '
    45 = '  +1.52%  '
    46 = '  +1.86%  '
    47 = '  +0.98%  '
    48 = '  +2.08%  '
    49 = '  +0.32%  '
    50 = '  +1.67%  '
    51 = '  +2.87%  '
}

foreach ($row in $priceUpdates.Keys) {
    $cell = $ws.Range("D$row")
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$row]
    $cell.ClearFormats()
}

foreach ($row in $volumeUpdates.Keys) {
    $ws.Range("E$row").Value = $volumeUpdates[$row]
}
